$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect so the cells below can be
# updated, then re-apply protection afterwards.
$ws.Unprotect()

# Update the disclaimer text date (2021-05-26 -> 2021-05-27) in cell A11
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."
# Re-fit the row height after the multi-line text write (avoid leaving a
# stray explicit row height behind).
$ws.Rows("11").AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.4984709084924442
$ws.Range("E2").Value = 0.004629940280480582

$ws.Range("D3").Value = 0.2464532146649872
$ws.Range("E3").Value = -0.003334299797042695

$ws.Range("D4").Value = 0.09538691754460515
$ws.Range("E4").Value = 0.005146228191289071

$ws.Range("D5").Value = 0.1024715561613304
$ws.Range("E5").Value = 0.00723293768545985

$ws.Range("D6").Value = 0.03005160826370766
$ws.Range("E6").Value = 0.01290077784101684

$ws.Range("D7").Value = 0.02716579487292538
$ws.Range("E7").Value = 0.007464409388226345

$ws.Range("E8").Value = 0.003308660594897539

# Restore sheet protection (matching the original protection settings as
# closely as the object model allows).
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
